$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 811, pushing existing rows 811:898 down to 813:900.
$ws.Rows("811:812").Insert()

# Populate the newly inserted row 811 with its data.
$ws.Cells.Item(811, 1).Value = 4
$ws.Cells.Item(811, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(811, 3).Value = "Los Lagos"
$ws.Cells.Item(811, 4).Value = 45212
$ws.Cells.Item(811, 5).Value = 10
$ws.Cells.Item(811, 6).Value = "Fruta"
$ws.Cells.Item(811, 7).Value = 100102
$ws.Cells.Item(811, 8).Value = "Cítricos"
$ws.Cells.Item(811, 9).Value = 100102005
$ws.Cells.Item(811, 10).Value = "Naranja"
$ws.Cells.Item(811, 11).Value = "Navel Late"
$ws.Cells.Item(811, 12).Value = "Primera"
$ws.Cells.Item(811, 13).Value = 300
$ws.Cells.Item(811, 14).Value = 18000
$ws.Cells.Item(811, 15).Value = 18000
$ws.Cells.Item(811, 16).Value = 18000
$ws.Cells.Item(811, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(811, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(811, 19).Value = 1200
$ws.Cells.Item(811, 20).Value = 15

# Populate the newly inserted row 812 with its data.
$ws.Cells.Item(812, 1).Value = 4
$ws.Cells.Item(812, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(812, 3).Value = "Los Lagos"
$ws.Cells.Item(812, 4).Value = 45212
$ws.Cells.Item(812, 5).Value = 10
$ws.Cells.Item(812, 6).Value = "Fruta"
$ws.Cells.Item(812, 7).Value = 100102
$ws.Cells.Item(812, 8).Value = "Cítricos"
$ws.Cells.Item(812, 9).Value = 100102005
$ws.Cells.Item(812, 10).Value = "Naranja"
$ws.Cells.Item(812, 11).Value = "Navel Late"
$ws.Cells.Item(812, 12).Value = "Segunda"
$ws.Cells.Item(812, 13).Value = 300
$ws.Cells.Item(812, 14).Value = 14000
$ws.Cells.Item(812, 15).Value = 14000
$ws.Cells.Item(812, 16).Value = 14000
$ws.Cells.Item(812, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(812, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(812, 19).Value = 933
$ws.Cells.Item(812, 20).Value = 15

# Column D holds a date/time serial value formatted as "YYYY-MM-DD HH:MM:SS" (style index 2),
# same as the rest of the column - match that number format on the two new rows.
$ws.Range("D811:D812").NumberFormat = "YYYY-MM-DD HH:MM:SS"
